$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10: 18 -> 1 (numeric)
$ws.Range("C10").Value = 1.0

# B11: "1" -> "R40"
$ws.Range("B11").Value = "R40"
